{"js": "// Replace the 25 two-digit-by-one-digit division prompts in the practice\n// table with the new set of problems, in document order. Each \"before\"\n// value is unique in the document, and the replacements are applied in\n// the same order they occur in the document, so a freshly-inserted value\n// (e.g. the new \"21\u00f79=\" introduced below) can never be picked up by a\n// later search for an older occurrence of the same text.\nconst replacements = [\n  [\"58\u00f77=\", \"57\u00f74=\"],\n  [\"21\u00f79=\", \"26\u00f78=\"],\n  [\"97\u00f74=\", \"55\u00f73=\"],\n  [\"11\u00f75=\", \"91\u00f75=\"],\n  [\"18\u00f74=\", \"36\u00f74=\"],\n  [\"35\u00f75=\", \"57\u00f76=\"],\n  [\"54\u00f75=\", \"51\u00f78=\"],\n  [\"26\u00f79=\", \"55\u00f78=\"],\n  [\"33\u00f74=\", \"94\u00f77=\"],\n  [\"10\u00f79=\", \"95\u00f77=\"],\n  [\"23\u00f76=\", \"90\u00f72=\"],\n  [\"72\u00f75=\", \"70\u00f79=\"],\n  [\"74\u00f79=\", \"67\u00f73=\"],\n  [\"37\u00f79=\", \"22\u00f72=\"],\n  [\"40\u00f76=\", \"21\u00f79=\"],\n  [\"45\u00f79=\", \"80\u00f72=\"],\n  [\"94\u00f72=\", \"52\u00f77=\"],\n  [\"19\u00f75=\", \"60\u00f74=\"],\n  [\"25\u00f77=\", \"13\u00f78=\"],\n  [\"22\u00f74=\", \"37\u00f75=\"],\n  [\"13\u00f79=\", \"57\u00f78=\"],\n  [\"15\u00f74=\", \"77\u00f75=\"],\n  [\"54\u00f72=\", \"25\u00f74=\"],\n  [\"75\u00f76=\", \"27\u00f78=\"],\n  [\"74\u00f76=\", \"71\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for '\" + before + \"'\");\n  }\n\n  // Each source string occurs exactly once in the document.\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-by-one-digit division prompts in the practice\n# table with the new set of problems, in document order. Each \"before\"\n# value is unique in the document, and the replacements are applied in\n# the same order they occur in the document, so a freshly-inserted value\n# (e.g. the new \"21\u00f79=\" introduced below) can never be picked up by a\n# later Find for an older occurrence of the same text.\n#\n# wdFindContinue = 1, wdReplaceOne = 1 (literal values; named constants\n# are not predefined in this host).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"58\u00f77=\", \"57\u00f74=\"),\n  @(\"21\u00f79=\", \"26\u00f78=\"),\n  @(\"97\u00f74=\", \"55\u00f73=\"),\n  @(\"11\u00f75=\", \"91\u00f75=\"),\n  @(\"18\u00f74=\", \"36\u00f74=\"),\n  @(\"35\u00f75=\", \"57\u00f76=\"),\n  @(\"54\u00f75=\", \"51\u00f78=\"),\n  @(\"26\u00f79=\", \"55\u00f78=\"),\n  @(\"33\u00f74=\", \"94\u00f77=\"),\n  @(\"10\u00f79=\", \"95\u00f77=\"),\n  @(\"23\u00f76=\", \"90\u00f72=\"),\n  @(\"72\u00f75=\", \"70\u00f79=\"),\n  @(\"74\u00f79=\", \"67\u00f73=\"),\n  @(\"37\u00f79=\", \"22\u00f72=\"),\n  @(\"40\u00f76=\", \"21\u00f79=\"),\n  @(\"45\u00f79=\", \"80\u00f72=\"),\n  @(\"94\u00f72=\", \"52\u00f77=\"),\n  @(\"19\u00f75=\", \"60\u00f74=\"),\n  @(\"25\u00f77=\", \"13\u00f78=\"),\n  @(\"22\u00f74=\", \"37\u00f75=\"),\n  @(\"13\u00f79=\", \"57\u00f78=\"),\n  @(\"15\u00f74=\", \"77\u00f75=\"),\n  @(\"54\u00f72=\", \"25\u00f74=\"),\n  @(\"75\u00f76=\", \"27\u00f78=\"),\n  @(\"74\u00f76=\", \"71\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n  $before = $pair[0]\n  $after = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($before, $false, $true, $false, $false, $false, $true, 1, $false, $after, 1)\n\n  if (-not $found) {\n    throw \"No match found for '$before'\"\n  }\n}\n"}
